# Layout + Routing v1 done
# Swap the Altimeter part (row 3) from the NXP MPL3115A2R1 to the
# STMicroelectronics LPS25HBTR, and add the new USB connector line (row 10).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: Altimeter part change (NXP MPL3115A2R1 -> ST LPS25HBTR) ---
$ws.Range("B3").Value = "LPS25HBTR"
$ws.Range("C3").Value = "STMicroelectronics"
$ws.Range("E3").Value = "HLGA-10L"
$ws.Range("F3").Value = "C87898"

# Give the MFR part # its own styling (matches the other "wrapped" part rows)
$ws.Range("B3").Style = "Normal"
$ws.Range("B3").Font.Name = "Arial"
$ws.Range("B3").Font.Size = 10

# F3 becomes a real hyperlink to the LCSC product page for the new part
$ws.Hyperlinks.Add(
    $ws.Range("F3"),
    "https://lcsc.com/product-detail/Pressure-Sensors_STMicroelectronics-LPS25HBTR_C87898.html/?href=jlc-SMT",
    [Type]::Missing,
    "https://lcsc.com/product-detail/Pressure-Sensors_STMicroelectronics-LPS25HBTR_C87898.html/?href=jlc-SMT"
) | Out-Null

# --- Row 10: new USB connector line item ---
$ws.Range("B10").Value = "1981568-1"
$ws.Range("C10").Value = "TE Connectivity"
$ws.Range("D10").Value = "USB"
$ws.Range("E10").Value = "USB Micro B"
$ws.Range("G10").Value = "USB Connector"

$ws.Rows.Item(10).RowHeight = 28.8

# H column holds a second, indented hyperlink copy of the designator
$ws.Columns.Item(8).ColumnWidth = 13.21875
$ws.Range("H10").Value = "1981568-1"
$ws.Range("H10").HorizontalAlignment = -4131
$ws.Range("H10").VerticalAlignment = -4108
$ws.Range("H10").WrapText = $true
$ws.Range("H10").IndentLevel = 1
$ws.Range("H10").Style = "Hyperlink"

$ws.Hyperlinks.Add(
    $ws.Range("H10"),
    "https://www.lcsc.com/product-detail/Molex_1981568-1_1981568-1.html"
) | Out-Null

$ws.Range("A10").Select()

$wb.Save()
